$wb = $excel.ActiveWorkbook

# --- 1. Rename Sheet1 -> Tests, add new "Tables" sheet right after it ---
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Tests"
$tbl = $wb.Worksheets.Add($null, $ws)
$tbl.Name = "Tables"

# --- 2. Seed new shared strings in the exact order needed so the
#        resulting sharedStrings.xml unique-string indices line up
#        with the target workbook (29..39). ---
$ws.Range("B24").Value = "Simulation"
$ws.Range("C24").Value = "PC support"
$ws.Range("B23").Value = "Simulation"
$ws.Range("C23").Value = "G Core support"
$ws.Range("E1").Value = "Testing  level"
$ws.Range("E18").Value = "Analysis"
$ws.Range("E2").Value = "Peer Review"
$ws.Range("E6").Value = "Demonstration"
$tbl.Range("A5").Value = "Quantitative"
$ws.Range("F1").Value = "Test reference"
$ws.Range("F6").Value = "T1001"
$ws.Range("F8").Value = "T1002"

# --- 3. Fill in the rest of the "Testing  level" (E) column on Tests ---
$ws.Range("E3").Value = "Peer Review"
$ws.Range("E5").Value = "Peer Review"
$ws.Range("E7").Value = "Peer Review"
$ws.Range("E8").Value = "Demonstration"
$ws.Range("E9").Value = "Demonstration"
$ws.Range("E10").Value = "Demonstration"
$ws.Range("E11").Value = "Demonstration"
$ws.Range("E12").Value = "Demonstration"
$ws.Range("E13").Value = "Peer Review"
$ws.Range("E14").Value = "Demonstration"
$ws.Range("E15").Value = "Peer Review"
$ws.Range("E16").Value = "Demonstration"
$ws.Range("E17").Value = "Demonstration"
$ws.Range("E19").Value = "Demonstration"
$ws.Range("E20").Value = "Demonstration"
$ws.Range("E21").Value = "Demonstration"
$ws.Range("E22").Value = "Demonstration"
$ws.Range("E23").Value = "Demonstration"
$ws.Range("E24").Value = "Demonstration"

# --- 4. New rows 23 & 24 on Tests (Requirement code / Topic / Subtopic) ---
$ws.Range("A23").Value = 1021
$ws.Range("A24").Value = 1022
# B23/C23 and B24/C24 already set above

# --- 5. Tables sheet lookup list (A2:A5) ---
$tbl.Range("A2").Value = "Peer Review"
$tbl.Range("A3").Value = "Analysis"
$tbl.Range("A4").Value = "Demonstration"
# A5 already set above

# --- 6. Column widths for the two new Tests columns, and for Tables!A ---
$ws.Columns.Item(5).ColumnWidth = 13.38
$ws.Columns.Item(6).ColumnWidth = 17.38
$tbl.Columns.Item(1).ColumnWidth = 10.36

# --- 7. Data validation: dropdown list sourced from Tables!$A$2:$A$5 ---
$dvRange = $ws.Range("E2:E1048576")
$dvRange.Validation.Add(3, 1, 1, "Tables!`$A`$2:`$A`$5")
$dvRange.Validation.IgnoreBlank = $true
$dvRange.Validation.InCellDropdown = $true
$dvRange.Validation.ShowInput = $true
$dvRange.Validation.ShowError = $true

# --- 8. Selections matching the target file ---
$ws.Range("F22").Select()
$tbl.Range("F6").Select()
$ws.Activate()
